$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Total fees paid to IFIs" (column F) values on rows 2-5
$ws.Range("F2").Value = 761579.37
$ws.Range("F3").Value = 761579.37
$ws.Range("F4").Value = 761579.37
$ws.Range("F5").Value = 761579.37

# Update selected cell in the sheet view from F8 to F2
$ws.Range("F2").Select()
